# "Minor adjustments + Bug fixes"
#  - insert a new "Dataset" column between "Batch Size" and "Command"
#    (the old "Command" column M shifts right to become column N)
#  - fill the new Dataset column for the existing rows with "CIFAR-10"
#  - add a new DenseNet-121 / ImageNet run in row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column M ("Command") one slot to the right, freeing up M for the
# new "Dataset" column. This naturally carries over formatting for every
# row (including the still-empty rows 5, 7-11).
$ws.Columns.Item(13).Insert(-4121)

# --- New "Dataset" column ---
$ws.Cells.Item(1, 13).Value = "Dataset"
$ws.Cells.Item(2, 13).Value = "CIFAR-10"
$ws.Cells.Item(3, 13).Value = "CIFAR-10"
$ws.Cells.Item(4, 13).Value = "CIFAR-10"

# --- New row 6: DenseNet-121 trained on ImageNet ---
# Start from row 4's formatting (closest fully-populated data row).
$ws.Range("A4:N4").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)
[void]$excel.CutCopyMode

$ws.Cells.Item(6, 13).Value = "ImageNet"
$ws.Cells.Item(6, 1).Value = "DenseNet-121"

$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Cells.Item(6, 4).Value = 0.0001
$ws.Cells.Item(6, 5).Value = 121
$ws.Cells.Item(6, 6).Value = 32
$ws.Cells.Item(6, 7).Value = 0

$ws.Range("H2").Copy()
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("H3").Copy()
$ws.Range("I6").PasteSpecial(-4163)
$ws.Range("J3").Copy()
$ws.Range("J6").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Cells.Item(6, 11).Value = 90
$ws.Cells.Item(6, 12).Value = 256
$ws.Cells.Item(6, 14).Value = "python main.py --layers 121 --growth 32 --reduce 0.5 --epochs 90 -b 256 --name DenseNet-121 --imagenet"

[void]$ws.Range("N6").Select()
